$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$footerArr = @($sec.Footers.Item(1))
Write-Host "arr count:" $footerArr.Count
